$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: new "Mussel, blue, raw" entry -------------------------------
# Text cells: set values first (apostrophe-prefix forces text for numeric-
# looking strings), then paste formats from existing cells that already
# carry the exact style combinations we need so no new fonts/fills are
# created and existing ones are reused by index.

# A15 - Source (no special style), reuse existing shared string.
$ws.Cells.Item(15, 1).Value2 = "Norwegian_Food_Composition_Table_2021"

# B15 - Food code "04.053" (looks numeric -> must force text), style
# matches J13 (fontId=1, fillId=2, centered).
$ws.Cells.Item(15, 2).Value2 = "'04.053"
$ws.Cells.Item(13, 10).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122)

# C15 - uk_name "Sea mussels", style matches C2 (quotePrefix, no fill).
$ws.Cells.Item(15, 3).Value2 = "Sea mussels"
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)

# D15 - food_name "Mussel, blue, raw", new style (fontId=1, fillId=2, no
# alignment) built by combining D13's font-only style with a yellow fill.
$ws.Cells.Item(15, 4).Value2 = "Mussel, blue, raw"
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4122)
$ws.Cells.Item(15, 4).Interior.Color = 65535

$excel.CutCopyMode = $false

# Numeric cells styled like J13 (fontId=1, fillId=2, centered).
$ws.Cells.Item(13, 10).Copy()
$styled14Cols = @(10, 12, 14, 16, 17, 18, 22, 23) # J, L, N, P, Q, R, V, W
foreach ($col in $styled14Cols) {
    $ws.Cells.Item(15, $col).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 10).Value2 = 14    # J15
$ws.Cells.Item(15, 11).Value2 = 0     # K15 (no style)
$ws.Cells.Item(15, 12).Value2 = 14    # L15
$ws.Cells.Item(15, 13).Value2 = 0     # M15 (no style)
$ws.Cells.Item(15, 14).Value2 = 0.8   # N15
$ws.Cells.Item(15, 16).Value2 = 0     # P15
$ws.Cells.Item(15, 17).Value2 = 0.27  # Q15
$ws.Cells.Item(15, 18).Value2 = 1.2   # R15
$ws.Cells.Item(15, 22).Value2 = 25    # V15
$ws.Cells.Item(15, 23).Value2 = 42    # W15
$ws.Cells.Item(15, 26).Value2 = 0     # Z15 (no style)

# --- Sheet view: dimension grows, selection moves to Z16 -----------------
$ws.Range("Z16").Select()
